# Apply the edits described in the commit:
# "New metadata from @wilsonsj100 and reworked steps"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the flux start time in B2 (13:37:30 -> 13:48:00)
$ws.Range("B2").Value = 0.57500000000000007

# Update L2's note from "SMARTX 440 use redo" to the new note "used redo "
$ws.Range("L2").Value = "used redo "

# Widen column C to fit the new note text
$ws.Range("C1").ColumnWidth = 15.875

# Move / update the active selection
$ws.Range("C12").Select()
